# Set the specific attendance-summary cells (columns D, E, G, H) to 1
# for the rows where they were previously 0, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    3  = @("G", "H")
    4  = @("H")
    5  = @("D", "E")
    6  = @("H")
    7  = @("H")
    8  = @("H")
    9  = @("D", "E")
    10 = @("H")
    11 = @("D", "E")
    12 = @("H")
    13 = @("H")
    14 = @("H")
    15 = @("H")
    16 = @("H")
    17 = @("H")
    18 = @("H")
}

foreach ($row in $changes.Keys) {
    foreach ($col in $changes[$row]) {
        $ws.Range("$col$row").Value = 1
    }
}
